$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'42.067.13"
$ws.Range("D2").NumberFormat = "General"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -1.43%  "
$ws.Range("E2").NumberFormat = "General"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'2.304.24"
$ws.Range("D3").NumberFormat = "General"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -2.13%  "
$ws.Range("E3").NumberFormat = "General"
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  -0.11%  "
$ws.Range("E4").NumberFormat = "General"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'318.04"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -0.05%  "
$ws.Range("E5").NumberFormat = "General"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'105.05"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -2.89%  "
$ws.Range("E6").NumberFormat = "General"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.630"
$ws.Range("D7").NumberFormat = "General"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -0.77%  "
$ws.Range("E7").NumberFormat = "General"
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'  -0.11%  "
$ws.Range("E8").NumberFormat = "General"
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'  -1.42%  "
$ws.Range("E9").NumberFormat = "General"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'39.84"
$ws.Range("D10").NumberFormat = "General"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -4.35%  "
$ws.Range("E10").NumberFormat = "General"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.0911"
$ws.Range("D11").NumberFormat = "General"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -1.59%  "
$ws.Range("E11").NumberFormat = "General"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'8.44"
$ws.Range("D12").NumberFormat = "General"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -0.63%  "
$ws.Range("E12").NumberFormat = "General"
$ws.Range("E12").Style = "Normal"
$ws.Range("E13").Value = "'  +0.41%  "
$ws.Range("E13").NumberFormat = "General"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.979"
$ws.Range("D14").NumberFormat = "General"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -1.89%  "
$ws.Range("E14").NumberFormat = "General"
$ws.Range("E14").Style = "Normal"
$ws.Range("E15").Value = "'  -2.97%  "
$ws.Range("E15").NumberFormat = "General"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'2.654.11"
$ws.Range("D16").NumberFormat = "General"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -2.07%  "
$ws.Range("E16").NumberFormat = "General"
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'2.302.42"
$ws.Range("D17").NumberFormat = "General"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -2.53%  "
$ws.Range("E17").NumberFormat = "General"
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'42.042.37"
$ws.Range("D18").NumberFormat = "General"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -1.50%  "
$ws.Range("E18").NumberFormat = "General"
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "'  +0.71%  "
$ws.Range("E19").NumberFormat = "General"
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = "'  -0.12%  "
$ws.Range("E20").NumberFormat = "General"
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'287.48"
$ws.Range("D21").NumberFormat = "General"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +11.94%  "
$ws.Range("E21").NumberFormat = "General"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'73.68"
$ws.Range("D22").NumberFormat = "General"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -3.77%  "
$ws.Range("E22").NumberFormat = "General"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'3.61"
$ws.Range("D23").NumberFormat = "General"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -0.28%  "
$ws.Range("E23").NumberFormat = "General"
$ws.Range("E23").Style = "Normal"
$ws.Range("B24").Value = "'InternetComputer(DFINITY)"
$ws.Range("B24").NumberFormat = "General"
$ws.Range("B24").Style = "Normal"
$ws.Range("C24").Value = "'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("C24").NumberFormat = "General"
$ws.Range("C24").Style = "Normal"
$ws.Range("D24").Value = "'10.14"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +8.16%  "
$ws.Range("E24").NumberFormat = "General"
$ws.Range("E24").Style = "Normal"
$ws.Range("B25").Value = "'ImmutableX"
$ws.Range("B25").NumberFormat = "General"
$ws.Range("B25").Style = "Normal"
$ws.Range("C25").Value = "'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("C25").NumberFormat = "General"
$ws.Range("C25").Style = "Normal"
$ws.Range("D25").Value = "'2.28"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -1.55%  "
$ws.Range("E25").NumberFormat = "General"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'1.00"
$ws.Range("D26").NumberFormat = "General"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -0.01%  "
$ws.Range("E26").NumberFormat = "General"
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'10.99"
$ws.Range("D27").NumberFormat = "General"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -3.86%  "
$ws.Range("E27").NumberFormat = "General"
$ws.Range("E27").Style = "Normal"
$ws.Range("B28").Value = "'EthereumClassic"
$ws.Range("B28").NumberFormat = "General"
$ws.Range("B28").Style = "Normal"
$ws.Range("C28").Value = "'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("C28").NumberFormat = "General"
$ws.Range("C28").Style = "Normal"
$ws.Range("D28").Value = "'23.47"
$ws.Range("D28").NumberFormat = "General"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +2.72%  "
$ws.Range("E28").NumberFormat = "General"
$ws.Range("E28").Style = "Normal"
$ws.Range("B29").Value = "'Toncoin"
$ws.Range("B29").NumberFormat = "General"
$ws.Range("B29").Style = "Normal"
$ws.Range("C29").Value = "'https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("C29").NumberFormat = "General"
$ws.Range("C29").Style = "Normal"
$ws.Range("D29").Value = "'2.23"
$ws.Range("D29").NumberFormat = "General"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -0.19%  "
$ws.Range("E29").NumberFormat = "General"
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'164.73"
$ws.Range("D30").NumberFormat = "General"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -5.62%  "
$ws.Range("E30").NumberFormat = "General"
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'35.44"
$ws.Range("D31").NumberFormat = "General"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -3.47%  "
$ws.Range("E31").NumberFormat = "General"
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'0.0884"
$ws.Range("D32").NumberFormat = "General"
$ws.Range("D32").Style = "Normal"
$ws.Range("E33").Value = "'  +0.29%  "
$ws.Range("E33").NumberFormat = "General"
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = "'  -3.20%  "
$ws.Range("E34").NumberFormat = "General"
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'0.133"
$ws.Range("D35").NumberFormat = "General"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +1.29%  "
$ws.Range("E35").NumberFormat = "General"
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'0.117"
$ws.Range("D36").NumberFormat = "General"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -8.65%  "
$ws.Range("E36").NumberFormat = "General"
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Value = "'  +0.70%  "
$ws.Range("E37").NumberFormat = "General"
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'2.95"
$ws.Range("D38").NumberFormat = "General"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +10.71%  "
$ws.Range("E38").NumberFormat = "General"
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'0.0353"
$ws.Range("D39").NumberFormat = "General"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -2.55%  "
$ws.Range("E39").NumberFormat = "General"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'3.62"
$ws.Range("D40").NumberFormat = "General"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -5.05%  "
$ws.Range("E40").NumberFormat = "General"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'102.41"
$ws.Range("D41").NumberFormat = "General"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +21.36%  "
$ws.Range("E41").NumberFormat = "General"
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = "'  +1.68%  "
$ws.Range("E42").NumberFormat = "General"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'71.08"
$ws.Range("D43").NumberFormat = "General"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -1.47%  "
$ws.Range("E43").NumberFormat = "General"
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = "'  -4.26%  "
$ws.Range("E44").NumberFormat = "General"
$ws.Range("E44").Style = "Normal"
$ws.Range("D46").Value = "'117.41"
$ws.Range("D46").NumberFormat = "General"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +3.50%  "
$ws.Range("E46").NumberFormat = "General"
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'12.13"
$ws.Range("D47").NumberFormat = "General"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +1.26%  "
$ws.Range("E47").NumberFormat = "General"
$ws.Range("E47").Style = "Normal"
$ws.Range("B48").Value = "'FraxShare"
$ws.Range("B48").NumberFormat = "General"
$ws.Range("B48").Style = "Normal"
$ws.Range("C48").Value = "'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("C48").NumberFormat = "General"
$ws.Range("C48").Style = "Normal"
$ws.Range("D48").Value = "'9.14"
$ws.Range("D48").NumberFormat = "General"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +0.13%  "
$ws.Range("E48").NumberFormat = "General"
$ws.Range("E48").Style = "Normal"
$ws.Range("B49").Value = "'ordi"
$ws.Range("B49").NumberFormat = "General"
$ws.Range("B49").Style = "Normal"
$ws.Range("C49").Value = "'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("C49").NumberFormat = "General"
$ws.Range("C49").Style = "Normal"
$ws.Range("D49").Value = "'77.82"
$ws.Range("D49").NumberFormat = "General"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +4.80%  "
$ws.Range("E49").NumberFormat = "General"
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = "'  -2.58%  "
$ws.Range("E50").NumberFormat = "General"
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = "'  +1.29%  "
$ws.Range("E51").NumberFormat = "General"
$ws.Range("E51").Style = "Normal"
